# Update the "dSF" column (F) with repulled/recalculated data.
# Diff shows only column F values changing for a subset of rows
# (the rest of the row data, including column E "dS0", stays the same).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = -4
    "F3"  = -3
    "F4"  = -3
    "F5"  = -11
    "F8"  = -2
    "F9"  = -3
    "F10" = -5
    "F14" = -3
    "F18" = -2
    "F22" = 9
    "F23" = -1
    "F26" = 0
    "F27" = 0
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
